$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old grid spanned A1:F9 (6 "mesa/taula" columns x up to 8 name rows).
# The new, trimmed-down test fixture only needs a 3x4 grid, so drop the
# now-unused rows (5-9) and columns (D-F) entirely first - this also gets
# rid of their stale column-width metadata automatically.
$ws.Range("A5:A9").EntireRow.Delete()
$ws.Range("D1:F9").EntireColumn.Delete()

# Rewrite the surviving cells with the updated names (some gain/lose a
# leading or trailing space, e.g. "Amaia" -> " Amaia"). Order matters here
# only insofar as it matches the order these strings were (re)typed by the
# original author, which is what's reflected in the saved file.
$ws.Range("C2").Value = "Jose "
$ws.Range("B2").Value = " Amaia"
$ws.Range("B3").Value = " Pepe - x - Celiaco"
$ws.Range("C3").Value = "Pepa - x "
$ws.Range("A4").Value = "Amaia-x-Celiaca "

# Row 4 used to continue into B4/C4 ("Pepe"); the new fixture leaves them
# blank.
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()

# Leave the selection where the author last left it when saving.
$ws.Range("C6").Select()
